# Anil changes after payment module testing.
# Clear out the "testing" values that were populated on Sheet1 (rows 2-4)
# for the EffectiveDate/ABGDate/PWBGPercentage/PWBGDate/Incoterms/Warranty/
# FEContent/TaxesAndDuties/GracePeriod columns.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

foreach ($r in 2..4) {
    # EffectiveDate (G), ABGDate (H), PWBGDate (J): keep the date formatting,
    # just blank out the values.
    $ws1.Range("G${r}:H${r}").ClearContents()
    $ws1.Range("J${r}").ClearContents()

    # PWBGPercentage (I), Incoterms (K), Warranty (L), FEContent (N),
    # TaxesAndDuties (O), GracePeriod (Q): remove the values AND the
    # formatting that was applied to them.
    $ws1.Range("I${r}").Clear()
    $ws1.Range("K${r}:L${r}").Clear()
    $ws1.Range("N${r}:O${r}").Clear()
    $ws1.Range("Q${r}").Clear()
}

# Reflect where the user was working when they made the change.
$ws1.Activate() | Out-Null
$ws1.Range("Q2:Q4").Select() | Out-Null
